# Update the "three-digit number divided by one-digit number" worksheet
# table by replacing each division expression with a new one, addressed by
# table row/column position (not by text) so that overlapping old/new
# values (e.g. "145÷4=" -> "133÷3=" and "133÷3=" -> "436÷5=") do not
# collide with each other.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row indices (1-based) that contain the division problems.
$rows = @(1, 5, 9, 13, 17)

$values = @(
    @("766÷8=", "589÷4=", "697÷3=", "180÷2=", "232÷6="),
    @("481÷7=", "831÷4=", "133÷3=", "175÷7=", "214÷2="),
    @("133÷9=", "229÷8=", "578÷7=", "255÷6=", "357÷3="),
    @("181÷7=", "957÷3=", "479÷7=", "113÷5=", "522÷8="),
    @("436÷5=", "720÷9=", "197÷9=", "695÷3=", "303÷5=")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowIndex = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $table.Cell($rowIndex, $c)
        $range = $cell.Range
        # Trim the trailing cell-mark / paragraph-mark characters so we
        # only touch the actual run text.
        $range.End = $range.End - 1
        $range.Text = $values[$r][$c - 1]
    }
}
